$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 381.42856
$ws.Range("I4").Value = 194
$ws.Range("J4").Value = 850
$ws.Range("K4").Value = 194
$ws.Range("L4").Value = 850
$ws.Range("M4").Value = -80
$ws.Range("N4").Value = -1078

$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()

$ws.Range("H64").Value = 4033.3333
$ws.Range("J64").Value = 4500
$ws.Range("L64").Value = 4500
$ws.Range("N64").Value = -4996

$ws.Range("H67").Value = 4033.3333
$ws.Range("J67").Value = 4500
$ws.Range("L67").Value = 4500
$ws.Range("N67").Value = -6216

$ws.Range("H74").Value = 5298
$ws.Range("I74").Value = 5000
$ws.Range("K74").Value = 5000
$ws.Range("M74").Value = -4064

$ws.Range("H77").Value = 5298
$ws.Range("I77").Value = 5000
$ws.Range("K77").Value = 25000
$ws.Range("M77").Value = -20320

$ws.Range("I80").Value = 617.1177
$ws.Range("K80").Value = 1851.3531
$ws.Range("M80").Value = -853.3531

$ws.Range("I83").Value = 617.1177
$ws.Range("K83").Value = 5554.0593
$ws.Range("M83").Value = -562.0592999999999

$ws.Range("H125").Value = 588.2857
$ws.Range("I125").Value = 533
$ws.Range("J125").Value = 662
$ws.Range("K125").Value = 4797
$ws.Range("L125").Value = 5958
$ws.Range("M125").Value = -2337
$ws.Range("N125").Value = -10878

$ws.Range("H138").Value = 2167.5696
$ws.Range("I138").Value = 1921.1875
$ws.Range("J138").Value = 2230.1428
$ws.Range("K138").Value = 5763.5625
$ws.Range("L138").Value = 6690.428400000001
$ws.Range("M138").Value = -623.5625
$ws.Range("N138").Value = -16970.4284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1900.6666
$ws.Range("I63").Value = 2080.8
$ws.Range("J63").Value = 1000
$ws.Range("K63").Value = 2080.8
$ws.Range("L63").Value = 1000
$ws.Range("M63").Value = -1394.8
$ws.Range("N63").Value = -2372

$ws.Range("H66").Value = 1900.6666
$ws.Range("I66").Value = 2080.8
$ws.Range("J66").Value = 1000
$ws.Range("K66").Value = 10404
$ws.Range("L66").Value = 5000
$ws.Range("M66").Value = -6972
$ws.Range("N66").Value = -11864

$ws.Range("H88").Value = 334384.66
$ws.Range("I88").Value = 1550
$ws.Range("K88").Value = 1550
$ws.Range("M88").Value = -1144

$ws.Range("H91").Value = 334384.66
$ws.Range("I91").Value = 1550
$ws.Range("K91").Value = 1550
$ws.Range("M91").Value = -146

$ws.Range("H97").Value = 1453.7273
$ws.Range("I97").Value = 1532.8235
$ws.Range("J97").Value = 1184.8
$ws.Range("K97").Value = 1532.8235
$ws.Range("L97").Value = 1184.8
$ws.Range("M97").Value = -1036.8235
$ws.Range("N97").Value = -2176.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 14660
$ws.Range("I26").Value = 6990
$ws.Range("J26").Value = 30000
$ws.Range("K26").Value = 6990
$ws.Range("L26").Value = 30000
$ws.Range("M26").Value = -6698
$ws.Range("N26").Value = -30584

$ws.Range("H86").Value = 1456.5938
$ws.Range("I86").Value = 1312.9584
$ws.Range("J86").Value = 1887.5
$ws.Range("K86").Value = 1312.9584
$ws.Range("L86").Value = 1887.5
$ws.Range("M86").Value = -189.9584
$ws.Range("N86").Value = -4133.5

$ws.Range("H89").Value = 1456.5938
$ws.Range("I89").Value = 1312.9584
$ws.Range("J89").Value = 1887.5
$ws.Range("K89").Value = 6564.791999999999
$ws.Range("L89").Value = 9437.5
$ws.Range("M89").Value = -948.7919999999995
$ws.Range("N89").Value = -20669.5

$ws.Range("H94").Value = 593.36365
$ws.Range("I94").Value = 593.36365
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 593.36365
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -142.36365
$ws.Range("N94").ClearContents()

$ws.Range("H134").Value = 4713.5356
$ws.Range("I134").Value = 4968.4585
$ws.Range("K134").Value = 14905.3755
$ws.Range("M134").Value = -12370.3755

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 26913.666
$ws.Range("J74").Value = 26913.666
$ws.Range("L74").Value = 26913.666
$ws.Range("N74").Value = -28661.666

$ws.Range("H77").Value = 26913.666
$ws.Range("J77").Value = 26913.666
$ws.Range("L77").Value = 80740.99800000001
$ws.Range("N77").Value = -89476.99800000001

$ws.Range("H132").Value = 18815.129
$ws.Range("I132").Value = 22725.541
$ws.Range("J132").Value = 5408
$ws.Range("K132").Value = 68176.62300000001
$ws.Range("L132").Value = 16224
$ws.Range("M132").Value = -65646.62300000001
$ws.Range("N132").Value = -21284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1455.5714
$ws.Range("I5").Value = 1455.5714
$ws.Range("K5").Value = 4366.7142
$ws.Range("M5").Value = -4254.7142

$ws.Range("H122").Value = 831.6
$ws.Range("J122").Value = 945.4375
$ws.Range("L122").Value = 8508.9375
$ws.Range("N122").Value = -13408.9375

$ws.Range("H131").Value = 186004.23
$ws.Range("I131").Value = 553.3333
$ws.Range("J131").Value = 196913.12
$ws.Range("K131").Value = 1659.9999
$ws.Range("L131").Value = 590739.36
$ws.Range("M131").Value = 3380.0001
$ws.Range("N131").Value = -600819.36

$ws.Range("H135").Value = 1455.5714
$ws.Range("I135").Value = 1455.5714
$ws.Range("K135").Value = 13100.1426
$ws.Range("M135").Value = -10565.1426

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3743.2856
$ws.Range("I80").Value = 3800
$ws.Range("J80").Value = 3727.818
$ws.Range("K80").Value = 3800
$ws.Range("L80").Value = 3727.818
$ws.Range("M80").Value = -2802
$ws.Range("N80").Value = -5723.818

$ws.Range("H83").Value = 3743.2856
$ws.Range("I83").Value = 3800
$ws.Range("J83").Value = 3727.818
$ws.Range("K83").Value = 19000
$ws.Range("L83").Value = 18639.09
$ws.Range("M83").Value = -14008
$ws.Range("N83").Value = -28623.09

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2355.4443
$ws.Range("I68").Value = 2385.7144
$ws.Range("K68").Value = 2385.7144
$ws.Range("M68").Value = -1636.7144

$ws.Range("H71").Value = 2355.4443
$ws.Range("I71").Value = 2385.7144
$ws.Range("K71").Value = 11928.572
$ws.Range("M71").Value = -8184.572

$ws.Range("H132").Value = 2061.3462
$ws.Range("I132").Value = 1849.8334
$ws.Range("J132").Value = 2537.25
$ws.Range("K132").Value = 5549.5002
$ws.Range("L132").Value = 7611.75
$ws.Range("M132").Value = -3019.5002
$ws.Range("N132").Value = -12671.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 166667460
$ws.Range("I81").Value = 946.2
$ws.Range("J81").Value = 1000000000
$ws.Range("K81").Value = 1892.4
$ws.Range("L81").Value = 2000000000
$ws.Range("M81").Value = -831.4000000000001
$ws.Range("N81").Value = -2000002122

$ws.Range("H84").Value = 166667460
$ws.Range("I84").Value = 946.2
$ws.Range("J84").Value = 1000000000
$ws.Range("K84").Value = 9462
$ws.Range("L84").Value = 10000000000
$ws.Range("M84").Value = -4158
$ws.Range("N84").Value = -10000010608
